# Updates cryptos list values per "Updated cryptos list" data refresh commit.
# D column holds prices as text (e.g. "25.944.83"); COM auto-converts numeric-looking
# strings to real numbers, so we prefix with a literal leading apostrophe to force text,
# exactly as Excel does when a user types a quoted numeric string into a General cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $ws.Range($range).Value = "'" + $text
}

# Row 2
Set-TextValue "D2" '25.944.83'
$ws.Range("E2").Value = '  +0.35%  '

# Row 3
Set-TextValue "D3" '1.735.43'
$ws.Range("E3").Value = '  -0.12%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
Set-TextValue "D5" '245.90'
$ws.Range("E5").Value = '  +3.66%  '

# Row 6
Set-TextValue "D6" '1.002'
$ws.Range("E6").Value = '  +0.17%  '

# Row 7
Set-TextValue "D7" '0.5013'
$ws.Range("E7").Value = '  -2.33%  '

# Row 8
Set-TextValue "D8" '0.2725'
$ws.Range("E8").Value = '  -0.38%  '

# Row 9
Set-TextValue "D9" '0.06182'
$ws.Range("E9").Value = '  +1.10%  '

# Row 10
Set-TextValue "D10" '1.747.65'
$ws.Range("E10").Value = '  +0.57%  '

# Row 11
Set-TextValue "D11" '0.07256'
$ws.Range("E11").Value = '  +1.34%  '

# Row 12
Set-TextValue "D12" '0.6532'
$ws.Range("E12").Value = '  +2.56%  '

# Row 13
Set-TextValue "D13" '15.10'
$ws.Range("E13").Value = '  +0.73%  '

# Row 14
Set-TextValue "D14" '4.736'
$ws.Range("E14").Value = '  +3.05%  '

# Row 15
Set-TextValue "D15" '77.56'
$ws.Range("E15").Value = '  +0.40%  '

# Row 16
Set-TextValue "D16" '1.001'
$ws.Range("E16").Value = '  +0.14%  '

# Row 17
Set-TextValue "D17" '1.002'
$ws.Range("E17").Value = '  +0.18%  '

# Row 18
Set-TextValue "D18" '25.963.09'
$ws.Range("E18").Value = '  +0.39%  '

# Row 19
Set-TextValue "D19" '11.88'
$ws.Range("E19").Value = '  +1.21%  '

# Row 20
Set-TextValue "D20" '0.000006813'
$ws.Range("E20").Value = '  +1.12%  '

# Row 21 — WrappedliquidstakedEther2.0
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D21" '1.966.41'
$ws.Range("E21").Value = '  +0.30%  '

# Row 22 — Uniswap
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D22" '4.591'
$ws.Range("E22").Value = '  +7.92%  '

# Row 23
Set-TextValue "D23" '8.773'
$ws.Range("E23").Value = '  +1.25%  '

# Row 24
Set-TextValue "D24" '5.399'
$ws.Range("E24").Value = '  +3.10%  '

# Row 25
Set-TextValue "D25" '133.72'
$ws.Range("E25").Value = '  -3.68%  '

# Row 26
Set-TextValue "D26" '1.507'
$ws.Range("E26").Value = '  -0.62%  '

# Row 27
Set-TextValue "D27" '15.27'
$ws.Range("E27").Value = '  +0.91%  '

# Row 28
Set-TextValue "D28" '1.780'
$ws.Range("E28").Value = '  +1.47%  '

# Row 29
Set-TextValue "D29" '105.63'
$ws.Range("E29").Value = '  +0.13%  '

# Row 30
Set-TextValue "D30" '3.974'
$ws.Range("E30").Value = '  -0.17%  '

# Row 31
Set-TextValue "D31" '0.08130'
$ws.Range("E31").Value = '  -2.46%  '

# Row 32
Set-TextValue "D32" '3.694'
$ws.Range("E32").Value = '  +1.27%  '

# Row 33
Set-TextValue "D33" '0.04731'
$ws.Range("E33").Value = '  +3.71%  '

# Row 34
Set-TextValue "D34" '2.664'
$ws.Range("E34").Value = '  +0.35%  '

# Row 35
Set-TextValue "D35" '0.9950'
$ws.Range("E35").Value = '  +1.04%  '

# Row 36
Set-TextValue "D36" '0.6075'
$ws.Range("E36").Value = '  -1.40%  '

# Row 37
Set-TextValue "D37" '2.733'
$ws.Range("E37").Value = '  +1.39%  '

# Row 38
Set-TextValue "D38" '0.01610'
$ws.Range("E38").Value = '  +1.07%  '

# Row 39
$ws.Range("E39").Value = '  +0.87%  '

# Row 40
Set-TextValue "D40" '1.001'
$ws.Range("E40").Value = '  +0.18%  '

# Row 41 — TrustWalletToken
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D41" '0.8231'
$ws.Range("E41").Value = '  +11.88%  '

# Row 42 — Quant
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D42" '100.66'
$ws.Range("E42").Value = '  +2.84%  '

# Row 43
Set-TextValue "D43" '0.3911'
$ws.Range("E43").Value = '  +1.75%  '

# Row 44
Set-TextValue "D44" '5.008'
$ws.Range("E44").Value = '  +1.13%  '

# Row 45
Set-TextValue "D45" '0.1172'
$ws.Range("E45").Value = '  +4.28%  '

# Row 46
Set-TextValue "D46" '6.341'
$ws.Range("E46").Value = '  +2.63%  '

# Row 47
Set-TextValue "D47" '55.66'
$ws.Range("E47").Value = '  +1.57%  '

# Row 48
Set-TextValue "D48" '0.05283'
$ws.Range("E48").Value = '  +0.30%  '

# Row 49
Set-TextValue "D49" '30.83'
$ws.Range("E49").Value = '  +0.98%  '

# Row 50
Set-TextValue "D50" '0.3466'
$ws.Range("E50").Value = '  +1.43%  '

# Row 51
Set-TextValue "D51" '7.595'
$ws.Range("E51").Value = '  +0.68%  '

